$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The PHYSICS and LINGUISTICS subject rows move up one row each,
# and MATHEMATICS moves down to where LINGUISTICS used to be.
# Resulting layout (rows 2-4), row 5 (MEDICINE) stays unchanged:
#   Row2: PHYSICS
#   Row3: LINGUISTICS
#   Row4: MATHEMATICS

$ws.Range("A2").Value = "PHYSICS"
$ws.Range("B2").Value = 36.29999923706055
$ws.Range("C2").Value = 8.0
$ws.Range("D2").Value = 2.0
$ws.Range("E2").Value = "Московский Выдуманный Университет, Московский Придуманный Институт"

$ws.Range("A3").Value = "LINGUISTICS"
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 0.0
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = "Воронежский Литературно-Переводческий Университет"

$ws.Range("A4").Value = "MATHEMATICS"
$ws.Range("B4").Value = 0.0
$ws.Range("C4").Value = 0.0
$ws.Range("D4").Value = 1.0
$ws.Range("E4").Value = "Казанский Университет Вычислений"
